$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 24 for the new "Chirnside Park" entry.
#    This shifts the existing rows 24-51 down to 25-52.
$ws.Rows.Item(24).Insert()
$ws.Range("A24").Value = "Chirnside Park"
$ws.Range("B24").Value = "Woolworths  239-241 Maroondah Hwy, Chirnside Park"
$ws.Range("C24").Value = "31/12/2020 10:00am - 10:15am"
$ws.Range("D24").Value = "Case shopped at venue"

# 2. Fix the Kmart site description on what is now row 25 (was row 24).
$ws.Range("B25").Value = "Kmart - 2107 Dandenong Road, Clayton"

# 3. Fix typo "vistied" -> "visited" on what is now row 27 (was row 26).
$ws.Range("D27").Value = "Case visited venue"

# 4. Tidy the exposure-period text on what is now row 30 (was row 29).
$ws.Range("C30").Value = "31/12/20 2pm-3pm"

# 5. Insert a new row at position 33 for the new "Keysborough" entry.
#    This shifts the current rows 33-52 down to 34-53.
$ws.Rows.Item(33).Insert()
$ws.Range("A33").Value = "Keysborough"
$ws.Range("B33").Value = "Sikh Temple Keysborough  200 Perry Road, Keysborough"
$ws.Range("C33").Value = "1/01/21 2:00pm-5:00pm"
$ws.Range("D33").Value = "Case visited venue"
